$ErrorActionPreference = "Stop"
$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 1: simplify the run structure of the
# "Run 02_prepare_to_check_Plutino_resonances.py on your local machine."
# paragraph so the bold filename is one run with a trailing space,
# merged with the leading "Run " run.
# ------------------------------------------------------------------
$targetPara1 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*prepare_to_check_Plutino_resonances.py*") {
        $targetPara1 = $p
        break
    }
}
if ($targetPara1 -eq $null) {
    throw "Could not find the '02_prepare_to_check_Plutino_resonances.py' paragraph"
}

$r1 = $d.Range($targetPara1.Range.Start, $targetPara1.Range.End - 1)
$ooxml1 = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Run </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">02_prepare_to_check_Plutino_resonances.py </w:t></w:r><w:r><w:t>on your local machine.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
$r1.InsertXML($ooxml1)

# ------------------------------------------------------------------
# Change 2: add a new paragraph after
# "Delete the entire folder from your local machine and re-download it from the cluster."
# describing examination of the plot_False/plot_True pdf files.
# ------------------------------------------------------------------
$targetPara2 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Delete the entire folder from your local machine and re-download it from the cluster*") {
        $targetPara2 = $p
        break
    }
}
if ($targetPara2 -eq $null) {
    throw "Could not find the 're-download it from the cluster' paragraph"
}

$r2 = $targetPara2.Range
$r2.Collapse(0)
$ooxml2 = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Examine all files </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>plot_False_3_2</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>_[</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>MPC designation].pdf</w:t></w:r><w:r><w:t xml:space="preserve"> and </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>plot_True_3_2_[MPC designation].pdf</w:t></w:r><w:r><w:t xml:space="preserve"> on your local machine. Note if any objects seem to be classified incorrectly. Enter the designations of the false negatives (objects classified False that should be classified True) on line 43 of </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>03_read_Plutino_resonances.py</w:t></w:r><w:r><w:t xml:space="preserve"> and enter the designations of the false positives (objects classified True that should be classified False) on line 44.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
$r2.InsertXML($ooxml2)
